{"js": "// Update the \"Th\u00f4ng tin th\u00e0nh vi\u00ean nh\u00f3m\" table:\n//  - resize the \u0110\u1ecba ch\u1ec9 / Email / S\u1ed1 \u0111i\u1ec7n tho\u1ea1i columns\n//  - fill in the \u0110\u1ecba ch\u1ec9 (and, where given, Email) values for each member\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// The member-info table is the 2nd table in the document (index 1).\nconst table = tables.items[1];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Resize the 2nd/3rd/4th columns (\u0110\u1ecba ch\u1ec9 / Email / S\u1ed1 \u0111i\u1ec7n tho\u1ea1i).\n// columnWidth is in points; 20 twips == 1 point.\nconst headerCells = rows.items[0].cells;\nheaderCells.load(\"items\");\nawait context.sync();\n\nheaderCells.items[1].columnWidth = 1984 / 20; // \u0110\u1ecba ch\u1ec9  -> 99.2pt  (1984 dxa)\nheaderCells.items[2].columnWidth = 3261 / 20; // Email    -> 163.05pt (3261 dxa)\nheaderCells.items[3].columnWidth = 1545 / 20; // S\u1ed1 \u0111i\u1ec7n tho\u1ea1i -> 77.25pt (1545 dxa)\nawait context.sync();\n\n// New \u0110\u1ecba ch\u1ec9 / Email values, keyed by row index (row 0 is the header).\nconst addressByRow = {\n  1: \"Ho\u00e0ng Li\u1ec7t, Ho\u00e0ng Mai, H\u00e0 N\u1ed9i\",\n  2: \"H\u00e0 \u0110\u00f4ng, H\u00e0 N\u1ed9i\",\n  3: \"Ba \u0110\u00ecnh, H\u00e0 N\u1ed9i\",\n  4: \"Hai B\u00e0 Tr\u01b0ng, H\u00e0 N\u1ed9i\",\n  5: \"Thanh Xu\u00e2n, H\u00e0 N\u1ed9i\",\n};\nconst emailByRow = {\n  1: \"dung.la187225@sis.hust.edu.vn\",\n  2: \"hung.nt187238@sis.hust.edu.vn\",\n};\n\nfor (let r = 1; r < rows.items.length; r++) {\n  const cells = rows.items[r].cells;\n  cells.load(\"items\");\n  await context.sync();\n\n  if (Object.prototype.hasOwnProperty.call(addressByRow, r)) {\n    cells.items[1].value = addressByRow[r];\n  }\n  if (Object.prototype.hasOwnProperty.call(emailByRow, r)) {\n    cells.items[2].value = emailByRow[r];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the \"Th\u00f4ng tin th\u00e0nh vi\u00ean nh\u00f3m\" table:\n#  - resize the \u0110\u1ecba ch\u1ec9 / Email / S\u1ed1 \u0111i\u1ec7n tho\u1ea1i columns\n#  - fill in the \u0110\u1ecba ch\u1ec9 (and, where given, Email) values for each member\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(2)\n\n# Resize the 2nd/3rd/4th columns (\u0110\u1ecba ch\u1ec9 / Email / S\u1ed1 \u0111i\u1ec7n tho\u1ea1i).\n# Column.Width is in points; 20 twips == 1 point.\n$tbl.Columns.Item(2).Width = 1984 / 20   # \u0110\u1ecba ch\u1ec9       -> 99.2pt  (1984 dxa)\n$tbl.Columns.Item(3).Width = 3261 / 20   # Email         -> 163.05pt (3261 dxa)\n$tbl.Columns.Item(4).Width = 1545 / 20   # S\u1ed1 \u0111i\u1ec7n tho\u1ea1i -> 77.25pt (1545 dxa)\n\n# New \u0110\u1ecba ch\u1ec9 / Email values, keyed by (1-based) table row.\n$addressByRow = @{\n    2 = \"Ho\u00e0ng Li\u1ec7t, Ho\u00e0ng Mai, H\u00e0 N\u1ed9i\"\n    3 = \"H\u00e0 \u0110\u00f4ng, H\u00e0 N\u1ed9i\"\n    4 = \"Ba \u0110\u00ecnh, H\u00e0 N\u1ed9i\"\n    5 = \"Hai B\u00e0 Tr\u01b0ng, H\u00e0 N\u1ed9i\"\n    6 = \"Thanh Xu\u00e2n, H\u00e0 N\u1ed9i\"\n}\n$emailByRow = @{\n    2 = \"dung.la187225@sis.hust.edu.vn\"\n    3 = \"hung.nt187238@sis.hust.edu.vn\"\n}\n\nfor ($r = 2; $r -le $tbl.Rows.Count; $r++) {\n    if ($addressByRow.ContainsKey($r)) {\n        $tbl.Cell($r, 2).Range.Text = $addressByRow[$r]\n    }\n    if ($emailByRow.ContainsKey($r)) {\n        $tbl.Cell($r, 3).Range.Text = $emailByRow[$r]\n    }\n}\n"}
